# Added instrument to cashflows (#204)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InvestorKyc")

# Update "Send Kyc Form To User *" (column P) from "No" to "Yes" for rows 2 and 3
$ws.Range("P2").Value = "Yes"
$ws.Range("P3").Value = "Yes"

# Reflect the active cell selection recorded in the saved workbook
$ws.Range("P4").Select()
